# -----------------------------------------------------------------------
# Update with Correct Forecast output
#
# - Renames Sheet1 -> "Sales vs PO"
# - Adds three new sheets: "Weekly Growth", "Volume Insights", "Prediction Info"
# - On "Sales vs PO": inserts a new "Order Week" column (C) holding the
#   original weekly date, shifts the dates in column A forward a week and
#   moves the old "PO_Requested_Qty" column to D, zeroing it out there
#   (the non-zero PO figures now live on the "Weekly Growth" sheet).
# - Populates "Weekly Growth" with the weeks that actually had PO activity
#   plus a week-over-week Growth% column.
# - Populates "Volume Insights" with summary stats over the PO quantities.
# - Populates "Prediction Info" with the predicted next week PO quantity.
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---- source data (as it existed on the original single sheet) ---------
$oldA = @(45320,45327,45334,45341,45348,45355,45362,45369,45376,45383,45390,45397,45404,45411,45418,45425,45432,45439,45446,45453,45460,45467,45474,45481,45488,45495,45502,45509,45516,45523,45530,45537,45544,45551,45558,45565,45572,45579,45586,45593,45600,45607,45614,45621,45628,45635,45642,45649)
$oldB = @(0,0,7,6,4,7,7,5,5,8,6,6,6,8,4,4,3,2,3,4,1,3,1,7,8,3,1,4,2,51,3,0,3,0,6,0,3,4,0,7,3,0,0,1,7,5,0,8)
$oldC = @(0,144,0,240,0,16,0,0,8,0,32,0,32,0,0,16,0,0,0,0,0,0,0,0,0,0,0,16,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)

$n = $oldA.Length

# ---- sheets -------------------------------------------------------------
$wsSales = $wb.ActiveSheet
$wsSales.Name = "Sales vs PO"

$wsGrowth = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsSales)
$wsGrowth.Name = "Weekly Growth"

$wsVolume = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsGrowth)
$wsVolume.Name = "Volume Insights"

$wsPred = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsVolume)
$wsPred.Name = "Prediction Info"

# =========================================================================
# Sheet 1: "Sales vs PO"
# Insert a new column C ("Order Week") before the old PO_Requested_Qty
# column (which becomes D). New column A = old A + 6 (one week later);
# new column C = old column A (the original order week); new column D is
# zeroed out everywhere (the non-zero figures move to "Weekly Growth").
# =========================================================================
$wsSales.Columns.Item(3).Insert()

$wsSales.Cells.Item(1,3).Value = "Order Week"

$wsSales.Range("C2:C$($n+1)").NumberFormat = "YYYY-MM-DD HH:MM:SS"

for ($i = 0; $i -lt $n; $i++) {
    $r = $i + 2
    $wsSales.Cells.Item($r,1).Value = $oldA[$i] + 6
    $wsSales.Cells.Item($r,2).Value = $oldB[$i]
    $wsSales.Cells.Item($r,3).Value = $oldA[$i]
    $wsSales.Cells.Item($r,4).Value = 0
}

# =========================================================================
# Sheet 2: "Weekly Growth"
# Only the weeks that had a non-zero PO_Requested_Qty, plus the
# week-over-week percentage growth of that quantity.
# =========================================================================
# Reuse the bold header formatting from "Sales vs PO" row 1 (same style
# used by every header cell in the source workbook).
$wsSales.Range("A1:C1").Copy($wsGrowth.Range("A1:C1"))

$wsGrowth.Cells.Item(1,1).Value = "ds"
$wsGrowth.Cells.Item(1,2).Value = "PO_Requested_Qty"
$wsGrowth.Cells.Item(1,3).Value = "Growth%"

$growthDs = @()
$growthQty = @()
for ($i = 0; $i -lt $n; $i++) {
    if ($oldC[$i] -ne 0) {
        $growthDs += $oldA[$i]
        $growthQty += $oldC[$i]
    }
}

$gn = $growthDs.Length
$wsGrowth.Range("A2:A$($gn+1)").NumberFormat = "YYYY-MM-DD HH:MM:SS"

for ($i = 0; $i -lt $gn; $i++) {
    $r = $i + 2
    $wsGrowth.Cells.Item($r,1).Value = $growthDs[$i]
    $wsGrowth.Cells.Item($r,2).Value = $growthQty[$i]
    if ($i -eq 0) {
        $wsGrowth.Cells.Item($r,3).Value = 0
    } else {
        $prev = $growthQty[$i-1]
        $cur = $growthQty[$i]
        $wsGrowth.Cells.Item($r,3).Value = ($cur - $prev) / $prev * 100
    }
}

# =========================================================================
# Sheet 3: "Volume Insights"
# =========================================================================
$wsSales.Range("A1:D1").Copy($wsVolume.Range("A1:D1"))

$wsVolume.Cells.Item(1,1).Value = "Total_PO_Quantity"
$wsVolume.Cells.Item(1,2).Value = "Average_PO_Quantity"
$wsVolume.Cells.Item(1,3).Value = "Max_PO_Quantity"
$wsVolume.Cells.Item(1,4).Value = "Min_PO_Quantity"

$total = 0
$maxQty = $growthQty[0]
$minQty = $growthQty[0]
for ($i = 0; $i -lt $gn; $i++) {
    $total += $growthQty[$i]
    if ($growthQty[$i] -gt $maxQty) { $maxQty = $growthQty[$i] }
    if ($growthQty[$i] -lt $minQty) { $minQty = $growthQty[$i] }
}
$avg = $total / $gn

$wsVolume.Cells.Item(2,1).Value = $total
$wsVolume.Cells.Item(2,2).Value = $avg
$wsVolume.Cells.Item(2,3).Value = $maxQty
$wsVolume.Cells.Item(2,4).Value = $minQty

# =========================================================================
# Sheet 4: "Prediction Info"
# =========================================================================
$wsSales.Range("A1").Copy($wsPred.Range("A1"))

$wsPred.Cells.Item(1,1).Value = "Predicted_Next_Week_PO_Quantity"
$wsPred.Cells.Item(2,1).Value = 0

# ---- leave the first sheet selected/active, matching the source file ---
$wsSales.Activate()

Write-Output "done"
